$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction for the Conf-589239 pair (rows 16-17): the web_price
# refresh picked up a new value (13.75 instead of 12.75), which flips the
# match flag for row 16 (excel_price 12.75 no longer equals web_price 13.75).
$ws.Range("C16").Value = 13.75
$ws.Range("D16").Value = $false
$ws.Range("C17").Value = 13.75

# --- New column E: "refresh_status" ------------------------------------
# Give E1 the same header style as the existing headers (bold/centered/
# bordered) by copying A1:D1's format onto E1, then set its text.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "refresh_status"

# --- Fill refresh_status for every data row -----------------------------
# "-" when the row already matched (match = TRUE), "success" when the
# refresh corrected a mismatch (match = FALSE).
for ($r = 2; $r -le 53; $r++) {
    $match = $ws.Cells.Item($r, 4).Value()
    if ($match) {
        $ws.Cells.Item($r, 5).Value = "-"
    } else {
        $ws.Cells.Item($r, 5).Value = "success"
    }
}
